$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post that used to occupy row 202 ("「内面が変われば周りが変わる」") was
# removed from the source content. Delete that entire row; Excel shifts all
# subsequent rows (203:293) up by one, which also updates the used range
# dimension from A1:C293 to A1:C292 automatically.
$ws.Rows.Item(202).Delete()
